$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy number formats from existing similarly-formatted cells to avoid creating
# duplicate style entries.
$ws.Range("B18").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E18").Copy() | Out-Null
$ws.Range("E20:E21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row 20 - entry No 11
$ws.Range("A20").Value2 = 11
$ws.Range("B20").Value2 = 44580
$ws.Range("C20").Value2 = "RPA RLOGIC"
$ws.Range("D20").Value2 = "1. Generated the Accounting Statements for the three centers and the same had been triggered success"
$ws.Range("E20").Value2 = 1
$ws.Range("F20").Value2 = "Completed"

# Row 21 - continuation of entry No 11
$ws.Range("D21").Value2 = "2. Generated the P&L reports of the Dec21 for the three centers and shared to Rahman san for the verification"
$ws.Range("E21").Value2 = 1
$ws.Range("F21").Value2 = "Completed"
